$d = $word.ActiveDocument

# Locate the two paragraphs to remove:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startIdx -eq -1 -and $t -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startIdx = $i
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1) {
    $startPara = $d.Paragraphs.Item($startIdx)
    $endPara = $d.Paragraphs.Item($endIdx)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
